$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Column widths: widen G, add new H, I, J columns ---
# (ColumnWidth is stored internally offset by 5/6 of a character from the
# serialized OOXML width, so subtract that back out to hit the target widths.)
$ws.Columns.Item(7).ColumnWidth = 54.1640625 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 49 - (5/6)
$ws.Columns.Item(9).ColumnWidth = 47.1640625 - (5/6)
$ws.Columns.Item(10).ColumnWidth = 50.6640625 - (5/6)

# --- Prepare the new H1/I1 header cells (reuse G1's existing header style) ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)

# --- Prepare body-style (wrap/vcenter) formatting for the new cells, copying from column F ---
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("J2").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("I3").PasteSpecial(-4122)

$ws.Range("F4").Copy()
$ws.Range("G4").PasteSpecial(-4122)

# --- Fill in the new ATOC / AOS / STQA mini course tables ---
# Column G: ATOC (Advanced Theory of Computation)
$ws.Range("G1").Value = "ATOC"
$ws.Range("G3").Value = "Topic: name- Week 1,lectures- 1 Lecture, duration- 01:20;`nVideo: link- https://drive.google.com/file/d/19dipmMmEiCQmlntPtJAtKF_DfoD_sufk/preview, name- ATC Week#1, duration- 00:00;`nSlides: slide- Course.pdf;`nSlides: slide- Lec1.pdf;"
$ws.Range("G4").Value = "Topic: name- Week 2,lectures- 1 Lecture, duration- 01:44;`nVideo: link- https://drive.google.com/file/d/1LLPqYvPSZruiDnWQfo4C6cxPAwNAJBAc/preview, name- ATC Week#2, duration- 01:44;`nSlides: slide- Lec1.pdf;"

# Column H: AOS (Advanced Operating System)
$ws.Range("H1").Value = "AOS"
$ws.Range("H3").Value = "Topic: name- Week 1,lectures- 1 Lecture, duration- 01:37;`nVideo: link- https://drive.google.com/file/d/1YOq3u1x-LfM-TkOvDmc8IIrJ3l2MCywy/preview, name- AOS Week#1part1.mp4, duration- 01:09;`nVideo: link- https://drive.google.com/file/d/1Nh8zm9SXHbo87Ge5tT8hdPpfrUb1AlaK/preview, name- AOS Week#1part2.mp4, duration- 00:28;"

# Column I: STQA (Software Testing and Quality Assurance)
$ws.Range("I1").Value = "STQA"
$ws.Range("I3").Value = "Topic: name- Week 1,lectures- 1 Lecture, duration- 00:00;`nSlides: slide- Chapter1.pptx;"

# Row 2: subject info for the three new courses
$ws.Range("G2").Value = "Subject: Advanced Theory of Computation;`nInstructor: Dr Ali Arshad;`nClassSenior: +92 313 7408286;`nNote: This website is going premium in 3 to 4 days. You can register yourself by contacting Qasim at Rs 350 / month or Rs 1400 / course. In this payment all courses will be open.;`nCreditHours: 3.0;"
$ws.Range("H2").Value = "Subject: Advanced Operating System;`nInstructor: Dr Nasir Mehmood;`nClassSenior: 0323213212;`nNote: This website is going premium in 3 to 4 days. You can register yourself by contacting Qasim at Rs 350 / month or Rs 1400 one time / course duration. In this payment all courses will be open.;`nCreditHours: 3.0;"
$ws.Range("I2").Value = "Subject: Software Testing and Quality Assurance;`nInstructor: Dr Abdul Hannan;`nClassSenior: N/A;`nNote: This website is going premium in 3 to 4 days. You can register yourself by contacting Qasim at Rs 350 / month or Rs 1400 one time / course duration. In this payment all courses will be open.;`nCreditHours: 3.0;"

# --- Row heights (auto-fit equivalents baked in by the authoring app) ---
$ws.Rows.Item(2).RowHeight = 136
$ws.Rows.Item(3).RowHeight = 409.6
$ws.Rows.Item(16).RowHeight = 409.5

# --- Sheet view: scroll to show the new columns, select I2 ---
$ws.Range("I2").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
